$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (row 3, 4, 5)
$ws.Range("A3").Value = 10100
$ws.Range("A4").Value = 10101

$ws.Range("B3").Value = "在放油螺栓处放置大号19号套筒"
$ws.Range("B4").Value = "使用大号棘轮扳手拆卸放油螺栓"
$ws.Range("B5").Value = "恭喜你完成训练"

$ws.Range("A5").Value = "Success0"

$ws.Range("C3").Value = "a"
$ws.Range("C4").Value = "b"
$ws.Range("C5").Value = "c"

# Column B width to fit content (bestFit, target ~29.875)
$ws.Columns.Item(2).ColumnWidth = 29.1

# Set selection to D5 as in diff
$ws.Range("D5").Select()
